# Update "想去人数" (want-to-go count) values in F column on the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 322
$ws1.Range("F5").Value = 2985
$ws1.Range("F6").Value = 2012
$ws1.Range("F7").Value = 389
$ws1.Range("F8").Value = 137
$ws1.Range("F9").Value = 1113
$ws1.Range("F11").Value = 660
$ws1.Range("F12").Value = 63

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 322
$ws4.Range("F5").Value = 2985
$ws4.Range("F6").Value = 2012
$ws4.Range("F7").Value = 389
$ws4.Range("F9").Value = 137
$ws4.Range("F10").Value = 1113
$ws4.Range("F12").Value = 661
$ws4.Range("F13").Value = 63
